# Update countries & provincias Spain
# - Refresh the "last updated" timestamp.
# - Update case numbers for several countries (India, Afganistan, Hungria,
#   Estado de Palestina, Maldivas, Georgia).
# - Lesoto's totals rose enough to move it up the ranking (from position
#   199 to 193, ahead of Polinesia Francesa/Gambia/Macao/Islas Turcas y
#   Caicos/San Martin (Parte Francesa)/Puerto Rico, each of which drops one
#   spot), and Dominica/Fiyi swap order (tied totals, no number changes
#   needed there).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Julio de 2020 a las 09:02"

# --- Straightforward numeric refreshes (country stays on its row) ---
# India (row 7)
$ws.Range("B7").Value = 674312
$ws.Range("C7").Value = 408
$ws.Range("E7").Value = 245950

# Afganistan (row 47)
$ws.Range("B47").Value = 32951
$ws.Range("C47").Value = 279
$ws.Range("D47").Value = 19366
$ws.Range("E47").Value = 12721
$ws.Range("G47").Value = 38
$ws.Range("H47").Value = 864

# Hungria (row 96)
$ws.Range("B96").Value = 4183
$ws.Range("C96").Value = 9
$ws.Range("D96").Value = 2811
$ws.Range("E96").Value = 783

# Estado de Palestina (row 98)
$ws.Range("E98").Value = 3357
$ws.Range("G98").Value = 2
$ws.Range("H98").Value = 15

# Maldivas (row 108)
$ws.Range("D108").Value = 2030
$ws.Range("E108").Value = 395

# Georgia (row 141)
$ws.Range("B141").Value = 951
$ws.Range("C141").Value = 3
$ws.Range("D141").Value = 828

# --- Lesoto climbs the ranking: rows 189-195 keep their row number but
#     the country name + stats of each shift down one slot to make room
#     for Lesoto's new, higher totals at row 189. ---

# Row 189: now Lesoto (new totals)
$ws.Range("A189").Value = "Lesoto"
$ws.Range("B189").Value = 63
$ws.Range("C189").Value = 28
$ws.Range("D189").Value = 11
$ws.Range("E189").Value = 52
$ws.Range("H189").Value = 0

# Row 190: now Polinesia Francesa (formerly row 189's numbers)
$ws.Range("A190").Value = "Polinesia Francesa"
$ws.Range("B190").Value = 62
$ws.Range("C190").Value = 0
$ws.Range("D190").Value = 60
$ws.Range("E190").Value = 2
$ws.Range("H190").Value = 0

# Row 191: now Gambia (formerly row 190's numbers)
$ws.Range("A191").Value = "Gambia"
$ws.Range("B191").Value = 57
$ws.Range("C191").Value = 0
$ws.Range("D191").Value = 27
$ws.Range("E191").Value = 28
$ws.Range("H191").Value = 2

# Row 192: now Macao (formerly row 191's numbers)
$ws.Range("A192").Value = "Macao"
$ws.Range("B192").Value = 46
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 45
$ws.Range("E192").Value = 1
$ws.Range("H192").Value = 0

# Row 193: now Islas Turcas y Caicos (formerly row 192's numbers)
$ws.Range("A193").Value = "Islas Turcas y Caicos"
$ws.Range("B193").Value = 45
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 11
$ws.Range("E193").Value = 32
$ws.Range("H193").Value = 2

# Row 194: now San Martin (Parte Francesa) (formerly row 193's numbers)
$ws.Range("A194").Value = "San Martin (Parte Francesa)"
$ws.Range("B194").Value = 43
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 37
$ws.Range("E194").Value = 3
$ws.Range("H194").Value = 3

# Row 195: now Puerto Rico (formerly row 194's numbers)
$ws.Range("A195").Value = "Puerto Rico"
$ws.Range("B195").Value = 39
$ws.Range("C195").Value = 0
$ws.Range("D195").Value = 1
$ws.Range("E195").Value = 36
$ws.Range("H195").Value = 2

# --- Dominica overtakes Fiyi (tied totals, so only the name order swaps;
#     row 204 (Laos) is unaffected, rows 205/206 swap names only). ---
$ws.Range("A205").Value = "Dominica"
$ws.Range("A206").Value = "Fiyi"
